$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 56 (pushes the footnote rows 56-59 down to 57-60)
$ws.Rows.Item(56).Insert()

# Populate the new row with the "Letter Packet" / 3639 service entry.
# A leading apostrophe forces text storage (shared-string) for values that
# would otherwise be auto-detected as numbers ("" and "3639").
$ws.Range("A56").Value = "'"
$ws.Range("B56").Value = "'"
$ws.Range("C56").Value = "Letter Packet"
$ws.Range("D56").Value = "'3639"
$ws.Range("E56").Value = "'3639"
$ws.Range("F56").Value = "B2X"
$ws.Range("G56").Value = "-"
$ws.Range("H56").Value = "Yes"
$ws.Range("I56").Value = "-"
$ws.Range("J56").Value = "Yes"
$ws.Range("K56").Value = "Yes"
$ws.Range("L56").Value = "-"
$ws.Range("M56").Value = "-"
$ws.Range("N56").Value = "-"
$ws.Range("O56").Value = "NO"
$ws.Range("P56").Value = "SE, DK, FI, AD, AE, AG, AI, AL, AM, AN, AO, AQ, AR, AS, AT, AU, AW, AX, AZ, BA, BB, BD, BE, BF, BG, BH, BI, BJ, BL, BM, BN, BO, BQ, BR, BS, BV, BW, BY, BZ, CA, CC, CD, CF, CG, CH, CI, CK, CL, CM, CN, CO, CR, CU, CV, CW, CX, CY, CZ, DE, DJ, DM, DO, DZ, EC, EE, EG, EH, ER, ES, ET, FJ, FK, FM, FO, FR, GA, GB, GD, GE, GF, GG, GH, GI, GL, GM, GN, GP, GQ, GR, GS, GT, GU, GW, GY, HK, HM, HN, HR, HT, HU, ID, IE, IM, IN, IO, IQ, IR, IS, IT, JE, JM, JO, JP, KE, KG, KH, KI, KM, KN, KP, KR, KW, KY, KZ, LA, LB, LC, LI, LK, LR, LS, LT, LU, LV, MA, MC, MD, ME, MF, MG, MH, MK, ML, MM, MN, MO, MP, MQ, MR, MS, MT, MU, MV, MW, MX, MY, MZ, NA, NC, NE, NF, NG, NI, NL, NP, NR, NU, NZ, OM, PA, PE, PF, PG, PH, PK, PL, PM, PN, PR, PS, PT, PW, PY, QA, RE, RO, RS, RU, RW, SA, SB, SC, SG, SH, SI, SJ, SK, SL, SM, SN, SO, SR, ST, SV, SX, SZ, TC, TD, TF, TG, TH, TJ, TK, TL, TM, TN, TO, TR, TT, TV, TW, TZ, UA, UG, UM, US, UY, UZ, VA, VC, VE, VG, VI, VN, VU, WF, WS, XK, YT, ZA, ZM, ZW"
$ws.Range("Q56").Value = "-"

# Copy the formatting (font/border/style) of the row above so the new row
# matches the rest of the table (style index 1) instead of the blank style
# the row-insert produced.
$ws.Range("A55:Q55").Copy()
$ws.Range("A56:Q56").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# The table's AutoFilter range/used range needs to grow by one row to cover
# the newly-inserted row (A1:P60 -> A1:P61).
$ws.AutoFilterMode = $false
$ws.Range("A1:P61").AutoFilter() | Out-Null

# Keep the hidden _FilterDatabase defined name in sync with the autofilter.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='Booking & SG API'!`$A`$1:`$P`$61"
    }
}
